$d = $word.ActiveDocument

# --- Step 1: fix the original (soon to be second) image paragraph: add w:lang and w:lastRenderedPageBreak ---
$pImg = $d.Paragraphs.Item(1)
$pImg.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" w:rsidR="001969A4" w:rsidRDefault="00E42B01"><w:r><w:rPr><w:noProof/><w:lang w:eastAsia="es-MX"/></w:rPr><w:lastRenderedPageBreak/><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0"><wp:extent cx="7077075" cy="5462267"/><wp:effectExtent l="0" t="0" r="0" b="5715"/><wp:docPr id="1" name="Imagen 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr id="0" name="Picture 1"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId4"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="7113286" cy="5490215"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>')

# --- Step 2: merge runs in the "Asunto: Festejo de fin de año" paragraph ---
$pAsunto = $d.Paragraphs.Item(2)
$pAsunto.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00E42B01" w:rsidRDefault="00E42B01" w:rsidP="00E42B01"><w:pPr><w:ind w:firstLine="708"/><w:rPr><w:sz w:val="72"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="72"/></w:rPr><w:t xml:space="preserve">Asunto: </w:t></w:r><w:r><w:rPr><w:sz w:val="72"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Festejo de fin de año </w:t></w:r></w:p>')

# --- Step 3: merge runs in the "Se le comunica a todo el PERSONAL que" paragraph ---
$pComunica = $d.Paragraphs.Item(4)
$pComunica.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00E42B01" w:rsidRDefault="00E42B01" w:rsidP="00E42B01"><w:pPr><w:ind w:firstLine="708"/><w:rPr><w:b/><w:sz w:val="48"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="48"/></w:rPr><w:t xml:space="preserve">Se le comunica a todo el   </w:t></w:r><w:r w:rsidRPr="00E42B01"><w:rPr><w:b/><w:sz w:val="48"/><w:u w:val="single"/></w:rPr><w:t>PERSONAL</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="48"/></w:rPr><w:t xml:space="preserve"> que</w:t></w:r></w:p>')

# --- Step 4: merge runs in the "La cita es el día ... 2019" paragraph ---
$pCita = $d.Paragraphs.Item(5)
$pCita.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00E42B01" w:rsidRDefault="00E42B01" w:rsidP="00E42B01"><w:pPr><w:ind w:firstLine="708"/><w:rPr><w:b/><w:sz w:val="48"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="48"/></w:rPr><w:t>La cita es el día   30 de noviembre del 2019</w:t></w:r></w:p>')

# --- Step 5: remove the _GoBack bookmark currently at the end of the 2019 block ---
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# --- Step 6: insert the whole new "2021" block (picture + 4 paragraphs + blank) before everything ---
$insertStart = $d.Paragraphs.Item(1).Range.Start
$insertRange = $d.Range($insertStart, $insertStart)
$insertRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:r><w:rPr><w:noProof/><w:lang w:eastAsia="es-MX"/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="72612839" wp14:editId="145114EF"><wp:extent cx="7077075" cy="5462267"/><wp:effectExtent l="0" t="0" r="0" b="5715"/><wp:docPr id="2" name="Imagen 2"/><wp:cNvGraphicFramePr><a:graphicFrameLocks noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr id="0" name="Picture 1"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId4"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="7113286" cy="5490215"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="708"/><w:rPr><w:sz w:val="72"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="72"/></w:rPr><w:t xml:space="preserve">Asunto: </w:t></w:r><w:r><w:rPr><w:sz w:val="72"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Festejo de fin de año </w:t></w:r></w:p><w:p/><w:p><w:pPr><w:ind w:firstLine="708"/><w:rPr><w:b/><w:sz w:val="48"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="48"/></w:rPr><w:t xml:space="preserve">Se le comunica a todo el   </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="48"/><w:u w:val="single"/></w:rPr><w:t>PERSONAL</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="48"/></w:rPr><w:t xml:space="preserve"> que</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="708"/><w:rPr><w:b/><w:sz w:val="48"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="48"/></w:rPr><w:t>La cita es el día   27 de noviembre del 2021</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="708"/><w:rPr><w:b/><w:sz w:val="48"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="48"/></w:rPr><w:t>a las   6</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="48"/></w:rPr><w:t>:pm   en punto</w:t></w:r></w:p><w:p/>')

# --- Step 7: add the _GoBack bookmark inside the new block, between "a las   6" and ":pm   en punto" ---
$findRange = $d.Content
$null = $findRange.Find.Execute("a las   6", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($findRange.End, $findRange.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output ("Done. ParaCount=" + $d.Paragraphs.Count)
